$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.426.91"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").Value = "2.427.03"
$ws.Range("E3").Value = "  -0.32%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.97"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.30"
$ws.Range("E6").Value = "  -1.12%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  -0.31%  "

$ws.Range("D9").Value = "2.422.76"
$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("E10").Value = "  -1.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.156"
$ws.Range("E11").Value = "  +0.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("E12").Value = "  -2.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("E13").Value = "  -1.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.65"
$ws.Range("E14").Value = "  -0.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000177"
$ws.Range("E15").Value = "  -2.32%  "

$ws.Range("D16").Value = "2.860.10"
$ws.Range("E16").Value = "  -0.48%  "

$ws.Range("D17").Value = "62.397.65"
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").Value = "2.425.11"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.17"
$ws.Range("E19").Value = "  -0.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.23"
$ws.Range("E20").Value = "  +3.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.95"
$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.16"
$ws.Range("E22").Value = "  -0.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.05"
$ws.Range("E23").Value = "  +11.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -4.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.78"
$ws.Range("E25").Value = "  -3.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "607.16"
$ws.Range("E26").Value = "  +0.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.72"
$ws.Range("E27").Value = "  +1.25%  "

$ws.Range("D28").Value = "0.0₃0997"
$ws.Range("E28").Value = "  -1.22%  "

$ws.Range("D29").Value = "2.580.08"
$ws.Range("E29").Value = "  +1.00%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.48"
$ws.Range("E31").Value = "  +1.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.10"
$ws.Range("E32").Value = "  -4.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.138"
$ws.Range("E34").Value = "  -3.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.09"
$ws.Range("E35").Value = "  +3.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.50"
$ws.Range("E36").Value = "  -0.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.997"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.375"
$ws.Range("E38").Value = "  -1.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.64"
$ws.Range("E39").Value = "  -0.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.30"
$ws.Range("E40").Value = "  -1.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "144.62"
$ws.Range("E41").Value = "  -2.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.77"
$ws.Range("E42").Value = "  -3.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.56"
$ws.Range("E43").Value = "  +1.07%  "

$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.79"
$ws.Range("E45").Value = "  +0.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "146.66"
$ws.Range("E46").Value = "  -1.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.72"
$ws.Range("E47").Value = "  +0.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.71"
$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0527"
$ws.Range("E49").Value = "  -1.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.593"
$ws.Range("E50").Value = "  -1.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0230"
$ws.Range("E51").Value = "  -0.86%  "
